$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.202.81"
$ws.Range("E2").Value = "  +5.58%  "
$ws.Range("D3").Value = "1.787.10"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.89"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2692"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06293"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("D10").Value = "1.785.22"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07048"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6284"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.663"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "79.98"
$ws.Range("D16").Value = "28.165.78"
$ws.Range("E16").Value = "  +6.22%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007243"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").Value = "2.010.78"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.552"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.763"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.258"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.69"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.78"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.856"
$ws.Range("E27").Value = "  +4.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.70"
$ws.Range("E28").Value = "  +3.03%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.177"
$ws.Range("E30").Value = "  +6.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08289"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.773"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04885"
$ws.Range("E33").Value = "  +8.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.083"
$ws.Range("E34").Value = "  +8.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6553"
$ws.Range("E35").Value = "  +5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.620"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9466"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.614"
$ws.Range("E38").Value = "  +8.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.068"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.921"
$ws.Range("E40").Value = "  +6.20%  "
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9998"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.90"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +3.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.200"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1215"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05447"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.040"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.296"
$ws.Range("E49").Value = "  +5.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.75"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("E51").Value = "  +2.45%  "